$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text so Excel does not auto-convert numeric-looking
# strings (e.g. "1.00", "0.110", "76.080.63") into numbers and strip formatting.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '76.080.63'
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '2.935.54'
$ws.Range("E3").Value = '  +4.29%  '

$ws.Range("E4:E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '202.81'
$ws.Range("E5").Value = '  +8.28%  '

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '597.54'
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("E7:E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("E8").Value = '  +0.50%  '

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.199'
$ws.Range("E9").Value = '  +4.13%  '

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '2.934.40'
$ws.Range("E10").Value = '  +4.30%  '

$ws.Range("E11:E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +16.82%  '

$ws.Range("E12:E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.76%  '

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '4.95'
$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '3.474.55'
$ws.Range("E14").Value = '  +4.16%  '

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '28.25'
$ws.Range("E15").Value = '  +4.92%  '

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '75.980.38'
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000191'
$ws.Range("E17").Value = '  +2.19%  '

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '2.918.90'
$ws.Range("E18").Value = '  +3.95%  '

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '13.24'
$ws.Range("E19").Value = '  +7.79%  '

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '8.98'
$ws.Range("E20").Value = '  -0.71%  '

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = '374.41'
$ws.Range("E21").Value = '  -0.78%  '

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '2.31'
$ws.Range("E22").Value = '  +1.93%  '

$ws.Range("E23:E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.42%  '

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = '71.74'
$ws.Range("E24").Value = '  +1.32%  '

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = '3.084.59'
$ws.Range("E26").Value = '  +4.57%  '

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = '4.29'
$ws.Range("E27").Value = '  +2.71%  '

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '9.76'
$ws.Range("E28").Value = '  -0.66%  '

$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000109'
$ws.Range("E29").Value = '  +4.66%  '

$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("E31:E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("B32:E32").NumberFormat = "@"
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '7.83'
$ws.Range("E32").Value = '  +1.79%  '

$ws.Range("B33:E33").NumberFormat = "@"
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '502.98'
$ws.Range("E33").Value = '  -2.51%  '

$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  +3.11%  '

$ws.Range("E35:E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("B36:E36").NumberFormat = "@"
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '20.34'
$ws.Range("E36").Value = '  +1.91%  '

$ws.Range("B37:E37").NumberFormat = "@"
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '163.71'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.110'
$ws.Range("E38").Value = '  +27.19%  '

$ws.Range("B39:E39").NumberFormat = "@"
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = '0.381'
$ws.Range("E39").Value = '  +11.53%  '

$ws.Range("B40:E40").NumberFormat = "@"
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").Value = '19.63'
$ws.Range("E40").Value = '  +1.31%  '

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.114'
$ws.Range("E41").Value = '  -4.35%  '

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '182.49'
$ws.Range("E42").Value = '  -3.36%  '

$ws.Range("E43:E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '5.03'
$ws.Range("E44").Value = '  +0.56%  '

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '1.67'
$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = '40.27'
$ws.Range("E46").Value = '  +0.50%  '

$ws.Range("E47:E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '2.36'
$ws.Range("E48").Value = '  +1.38%  '

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = '0.581'
$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '3.76'
$ws.Range("E50").Value = '  +0.77%  '

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = '22.55'
$ws.Range("E51").Value = '  +7.75%  '
